# Update cryptocurrency Price (D) and Volume(1h) (E) columns to refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.114.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "'1.904.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'0.8382"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'242.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'0.9992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.3291"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'26.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").Value = "'0.07078"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "'0.7650"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "'1.908.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "'5.281"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "'92.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "'30.111.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "'5.886"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "'244.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").Value = "'0.000007777"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'2.160.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'0.9997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'7.017"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").Value = "'0.1773"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +25.36%  "
$ws.Range("D26").Value = "'9.306"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").Value = "'165.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").Value = "'18.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'2.102"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'1.365"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").Value = "'1.521"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'0.05950"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.55%  "
$ws.Range("D33").Value = "'4.301"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'4.086"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "'0.7342"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "'2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").Value = "'0.01924"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").Value = "'2.787"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "'0.4456"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'73.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "'5.960"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'0.9987"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "'102.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").Value = "'7.583"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "'9.835"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").Value = "'1.003.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").Value = "'2.056.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'1.520"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
